# Power Models and Textures.xlsx
# Commit: "Add information for waterspill_splash.igb"
#
# The table on Sheet1 is sorted alphabetically by column A (File). A new
# file, "textures\waterspill_splash.igb", needs to be inserted in its
# correctly-sorted position -- alphabetically just before
# "textures\white.igb", which lived at row 139 -- pushing the previous
# rows 139-143 down to 140-144.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Insert a new blank row at row 139 (shifts old rows 139:143 -> 140:144)
$ws.Rows.Item(139).Insert()

# Grow the table/autofilter range so the new row (and the now-last row,
# 144) are included as part of Table2.
$lo.Resize($ws.Range("A1:G144"))

# Fill in the new row's data -- same PC/GameCube/PS2/PSP/Xbox "source"
# columns used by the other stock-permanent-package rows, filed under
# note category 1b ("Used in all versions of XML2 but not originally in
# permanent").
$ws.Range("A139").Value = "textures\waterspill_splash.igb"
$ws.Range("B139").Value = "XML2 PC"
$ws.Range("C139").Value = "XML2 GameCube"
$ws.Range("D139").Value = "XML2 PS2"
$ws.Range("E139").Value = "XML2 PSP"
$ws.Range("F139").Value = "XML2 Xbox"
$ws.Range("G139").Value = "1b. Used in all versions of XML2 but not originally in permanent"

# The sheet's standalone conditional-formatting rules (outside the table
# definition itself) cover A2:A143 / B2:B143 / ... and need to grow to
# row 144 too, same as the table did. Walk the sheet-wide, index-stable
# FormatConditions collection (rather than re-querying per-column ranges,
# whose .FormatConditions membership shifts as earlier rules are resized)
# so each of the 13 existing rules is resized exactly once, in place.
$allConditions = $ws.Cells.FormatConditions
$newAppliesTo = @(
    "A2:A144", "A2:A144",
    "B2:B144", "B2:B144",
    "B2:F144",
    "C2:C144", "C2:C144",
    "D2:D144", "D2:D144",
    "E2:E144", "E2:E144",
    "F2:F144", "F2:F144"
)
for ($j = 1; $j -le $allConditions.Count; $j++) {
    $allConditions.Item($j).ModifyAppliesToRange($ws.Range($newAppliesTo[$j - 1]))
}

# Leave the cursor where the author apparently left it when saving.
$ws.Range("B10").Select()
